$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: add E11 ---
$ws.Range("E11").Value = "Refactored Host App"

# --- Row 12: add E12 ---
$ws.Range("E12").Value = "Wrote Two Extension"

# --- Row 13: new row (B13, C13, E13) ---
$ws.Range("B13").Value = "Dale "
$ws.Range("C13").Value = "Coding for Sprint 2"
$ws.Range("C13").VerticalAlignment = -4108
$ws.Range("E13").Value = "Hooked up extensions and host app"

# --- Row 16: new row (B16, C16) ---
$ws.Range("B16").Value = "Dale"
$ws.Range("C16").Value = "Module Decomp Diagram 2.0"
$ws.Range("C16").VerticalAlignment = -4108

# --- Row 17: new row (B17, C17) ---
$ws.Range("B17").Value = "Dale"
$ws.Range("C17").Value = "Dependency Diagram 2.0"
$ws.Range("C17").VerticalAlignment = -4108

# --- Row 18: new row (B18, C18) ---
$ws.Range("B18").Value = "Dale"
$ws.Range("C18").Value = "class Diagram 2.0"
$ws.Range("C18").VerticalAlignment = -4108

# --- Row 19: new row (B19, C19) ---
$ws.Range("B19").Value = "Dale"
$ws.Range("C19").Value = "Startup 2.0"
$ws.Range("C19").VerticalAlignment = -4108

# --- Update selection to match final state ---
$ws.Range("C19").Select()
